$d = $word.ActiveDocument

# Merge the two runs ("We should take care of " + "environment lol")
# back into the original single-run text "We should take care of enviroment".
$d.Content.Find.Execute("We should take care of environment lol", $false, $false, $false, $false, $false, $true, 1, $false, "We should take care of enviroment", 2)
